# Update "Fruta, Agrícola del Norte S.A. de Arica - Tuna" sheet with the
# latest weekly price data. Existing rows 4-13 are re-shuffled/updated and a
# new row (14) is appended, matching the target diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D,L,M,N,O,P,Q,R,S,T for rows 4..14 in the final/target layout.
$data = @(
    @{ Row = 4;  D = 44972; L = "Segunda"; M = 140; N = 27000; O = 28000; P = 27429; Q = "`$/caja 18 kilos"; R = "Región Metropolitana"; S = 1524; T = 18 },
    @{ Row = 5;  D = 44979; L = "Segunda"; M = 250; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 6;  D = 44664; L = "Segunda"; M = 150; N = 29000; O = 30000; P = 29500; Q = "`$/caja 18 kilos"; R = "Región de Coquimbo";   S = 1639; T = 18 },
    @{ Row = 7;  D = 44643; L = "Primera"; M = 160; N = 28000; O = 30000; P = 29000; Q = "`$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1450; T = 20 },
    @{ Row = 8;  D = 44671; L = "Segunda"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 9;  D = 44965; L = "Primera"; M = 100; N = 34000; O = 35000; P = 34600; Q = "`$/caja 18 kilos"; R = "Región de Coquimbo";   S = 1922; T = 18 },
    @{ Row = 10; D = 44965; L = "Segunda"; M = 120; N = 32000; O = 33000; P = 32333; Q = "`$/caja 18 kilos"; R = "Región de Coquimbo";   S = 1796; T = 18 },
    @{ Row = 11; D = 44636; L = "Primera"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 12; D = 44993; L = "Segunda"; M = 130; N = 25000; O = 26000; P = 25462; Q = "`$/caja 18 kilos"; R = "Región de Coquimbo";   S = 1273; T = 20 },
    @{ Row = 13; D = 44679; L = "Segunda"; M = 200; N = 29000; O = 30000; P = 29500; Q = "`$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1475; T = 20 },
    @{ Row = 14; D = 44679; L = "Tercera"; M = 200; N = 24000; O = 25000; P = 24500; Q = "`$/caja 20 kilos"; R = "Región de Coquimbo";   S = 1225; T = 20 }
)

# Row 12 is brand new: duplicate the static columns (A,B,C,E,F,G,H,I,J,K) from
# row 11 into the new row 14 (shifting rows 12-13 down by one first).
$ws.Rows.Item(12).Insert()

# Copy the repeating, non-changing columns from row 11 into the freshly
# inserted row 12 (A, B, C, E, F, G, H, I, J, K are constant for every data row).
$staticCols = @(1, 2, 3, 5, 6, 7, 8, 9, 10, 11)
foreach ($col in $staticCols) {
    $ws.Cells.Item(12, $col).Value = $ws.Cells.Item(11, $col).Value2
}
$ws.Cells.Item(12, 4).NumberFormat = $ws.Cells.Item(11, 4).NumberFormat

# Now write the D, L, M, N, O, P, Q, R, S, T values for every row 4..14.
foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Cells.Item($r, 4).Value = $entry.D
    $ws.Cells.Item($r, 12).Value = $entry.L
    $ws.Cells.Item($r, 13).Value = $entry.M
    $ws.Cells.Item($r, 14).Value = $entry.N
    $ws.Cells.Item($r, 15).Value = $entry.O
    $ws.Cells.Item($r, 16).Value = $entry.P
    $ws.Cells.Item($r, 17).Value = $entry.Q
    $ws.Cells.Item($r, 18).Value = $entry.R
    $ws.Cells.Item($r, 19).Value = $entry.S
    $ws.Cells.Item($r, 20).Value = $entry.T
}
